# Capital-Expenditure.xlsx update: WIFI modules update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2016")

# Remove the stray blank row 9 (between the blank spacer row 8 and the Total row) so
# everything below shifts up by one (Total row -> row 9, Grand total row -> row 10).
$ws.Rows(9).Delete()

# Row 3: fix the description text and add a detail description for the sensor row
$ws.Range("D3").Value = "Sensor: cảm biến ánh sáng, nhiệt độ, độ ẩm đất"
$ws.Range("H3").Value = "35k / 1 sensor"

# Row 5: add detail description text for "Lam nha kinh"
$ws.Range("H5").Value = "4 mat kinh mica (), 1 mat kinh mica ()"

# Row 6: was "Vavle dien tu" -> now "Module WIFI ESP8266" with a date and quantities
$ws.Range("C6").NumberFormat = "MM/DD/YY"
$ws.Range("C6").Value = "05/25/16"
$ws.Range("D6").Value = "Module WIFI ESP8266"
$ws.Range("E6").Value = 90
$ws.Range("F6").NumberFormat = "#,##0.000\ [$₫-42A];[RED]\-#,##0.000\ [$₫-42A]"
$ws.Range("F6").Value = 90
$ws.Range("G6").NumberFormat = "#,##0.000\ [$₫-42A];[RED]\-#,##0.000\ [$₫-42A]"
$ws.Range("G6").Value = 90

# Row 7: was "Module wifi ESP8266" -> now "AMS1117" with a date, quantities, and detail
$ws.Range("C7").NumberFormat = "MM/DD/YY"
$ws.Range("C7").Value = "05/25/16"
$ws.Range("D7").Value = "AMS1117"
$ws.Range("E7").Value = 10
$ws.Range("F7").NumberFormat = "#,##0.000\ [$₫-42A];[RED]\-#,##0.000\ [$₫-42A]"
$ws.Range("F7").Value = 10
$ws.Range("G7").NumberFormat = "#,##0.000\ [$₫-42A];[RED]\-#,##0.000\ [$₫-42A]"

$ws.Range("H7").Value = "Mach giam ap cho ESP 8266"

# Row 8 loses its stray formatting remnants (clear now-unused cells)
$ws.Range("C8").ClearContents()
$ws.Range("E8").NumberFormat = "General"
$ws.Range("H8").ClearContents()

$wb.Save()
